$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header timestamp text
$ws.Range("A1").Value = "Datos actualizados a 20 de Abril de 2020 a las 03:52"

# Row 4 - Estados Unidos: numeric refresh
$ws.Range("B4").Value = 763836
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 71012
$ws.Range("E4").Value = 652269
$ws.Range("F4").Value = 13566
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 40555

# Rows 116-120: Guatemala moves up (between Isla de Man and Sri Lanka); everything below
# shifts down by one row through Vietnam, which now precedes Venezuela again.
$ws.Range("A116").Value = "Guatemala"
$ws.Range("B116").Value = 289
$ws.Range("C116").Value = 32
$ws.Range("D116").Value = 21
$ws.Range("E116").Value = 261
$ws.Range("F116").Value = 3
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 7

$ws.Range("A117").Value = "Sri Lanka"
$ws.Range("B117").Value = 271
$ws.Range("C117").Value = 0
$ws.Range("D117").Value = 96
$ws.Range("E117").Value = 168
$ws.Range("F117").Value = 1
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 7

$ws.Range("A118").Value = "Mayotte"
$ws.Range("B118").Value = 271
$ws.Range("C118").Value = 0
$ws.Range("D118").Value = 117
$ws.Range("E118").Value = 150
$ws.Range("F118").Value = 5
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 4

$ws.Range("A119").Value = "Kenia"
$ws.Range("B119").Value = 270
$ws.Range("C119").Value = 0
$ws.Range("D119").Value = 67
$ws.Range("E119").Value = 189
$ws.Range("F119").Value = 2
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 14

$ws.Range("A120").Value = "Vietnam"
$ws.Range("B120").Value = 268
$ws.Range("C120").Value = 0
$ws.Range("D120").Value = 202
$ws.Range("E120").Value = 66
$ws.Range("F120").Value = 8
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 0

# Row 121 (Venezuela) is unchanged.

# Rows 189-190: San Cristobal y Nieves moves ahead of Santa Lucia.
$ws.Range("A189").Value = "San Cristobal y Nieves"
$ws.Range("B189").Value = 15
$ws.Range("C189").Value = 1
$ws.Range("D189").Value = 0
$ws.Range("E189").Value = 15
$ws.Range("F189").Value = 0
$ws.Range("G189").Value = 0
$ws.Range("H189").Value = 0

$ws.Range("A190").Value = "Santa Lucia"
$ws.Range("B190").Value = 15
$ws.Range("C190").Value = 0
$ws.Range("D190").Value = 11
$ws.Range("E190").Value = 4
$ws.Range("F190").Value = 0
$ws.Range("G190").Value = 0
$ws.Range("H190").Value = 0
